$d = $word.ActiveDocument

# 1) Typo fix: capitalization of "Lettre Recommandée avec Accusé de Réception"
#    -> "Lettre recommandée avec accusé de réception"
$d.Content.Find.Execute("Lettre Recommand", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lettre recommand", 2)
$d.Content.Find.Execute("e avec Accus", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e avec accus", 2)
$d.Content.Find.Execute("de R" + [char]0x00E9 + "ception", $true, $false, $false, $false, $false,
                         $true, 1, $false, "de r" + [char]0x00E9 + "ception", 2)

# 2) Replace straight apostrophes with curly ones in the two remaining occurrences
$d.Content.Find.Execute("l'employeur", $true, $false, $false, $false, $false,
                         $true, 1, $false, "l" + [char]0x2019 + "employeur", 2)
$d.Content.Find.Execute("l'expression", $true, $false, $false, $false, $false,
                         $true, 1, $false, "l" + [char]0x2019 + "expression", 2)
